$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Set 1")

# --- Data corrections (DBCP datasource closed twice -> off-by-one counters) ---
$ws.Range("B5").Value  = 28672
$ws.Range("B6").Value  = 1401640

$ws.Range("B14").Value = 19999
$ws.Range("F14").Value = 19999

$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 2

$ws.Range("B16").Value = 1220000
$ws.Range("F16").Value = 1220000

$ws.Range("C17").Value = 12

# --- View / selection state ---
$ws.Activate()
$ws.Range("C21").Select()
